# Update countries & provincias Spain
# - Update Bolivia's COVID stats
# - Update Honduras' COVID stats
# - Re-sort the table descending by "Casos totales" (column B)
# - Update the "Datos actualizados" timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the rows for Bolivia and Honduras in the country table (A4:A216)
$boliviaRow = $null
$hondurasRow = $null
for ($r = 4; $r -le 216; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($name -eq "Bolivia") { $boliviaRow = $r }
    if ($name -eq "Honduras") { $hondurasRow = $r }
    if ($boliviaRow -and $hondurasRow) { break }
}

# Update Bolivia: Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$ws.Cells.Item($boliviaRow, 2).Value2 = 1053
$ws.Cells.Item($boliviaRow, 3).Value2 = 39
$ws.Cells.Item($boliviaRow, 4).Value2 = 110
$ws.Cells.Item($boliviaRow, 5).Value2 = 888
$ws.Cells.Item($boliviaRow, 6).Value2 = 3
$ws.Cells.Item($boliviaRow, 7).Value2 = 2
$ws.Cells.Item($boliviaRow, 8).Value2 = 55

# Update Honduras: Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$ws.Cells.Item($hondurasRow, 2).Value2 = 738
$ws.Cells.Item($hondurasRow, 3).Value2 = 36
$ws.Cells.Item($hondurasRow, 4).Value2 = 79
$ws.Cells.Item($hondurasRow, 5).Value2 = 593
$ws.Cells.Item($hondurasRow, 6).Value2 = 10
$ws.Cells.Item($hondurasRow, 7).Value2 = 2
$ws.Cells.Item($hondurasRow, 8).Value2 = 66

# Re-sort the country table (A3:H216, with header row 3) descending by column B (Casos totales)
$rng = $ws.Range("A3:H216")
$rng.Sort($ws.Range("B3"), 2, $null, $null, 1, $null, 1, 1)

# Update the "updated at" timestamp banner in A1
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 29 de Abril de 2020 a las 04:52"
